$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

# Rename remaining sheets (case/diacritics updates)
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the dashboard sheet as the active/selected tab (it was the
# active sheet before the edit, and deleting the last sheet would
# otherwise shift selection onto the new last tab)
$wb.Worksheets.Item("PAINEIS DARQ").Activate() | Out-Null
